$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 3 and row 4:
# Row 3 becomes what row 4 used to be, and row 4 becomes what row 3 used to be.

# Date (column D)
$ws.Range("D3").Value2 = 44672
$ws.Range("D4").Value2 = 44993

# Volumen (column M)
$ws.Range("M3").Value2 = 8
$ws.Range("M4").Value2 = 14

# Precio máximo (column O)
$ws.Range("O3").Value2 = 180000
$ws.Range("O4").Value2 = 200000

# Precio promedio ponderado (column P)
$ws.Range("P3").Value2 = 180000
$ws.Range("P4").Value2 = 190000

# Precio $/Kg (column S)
$ws.Range("S3").Value2 = 180000
$ws.Range("S4").Value2 = 190000
